$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix Q calculation: remove the erroneous *100 multiplication (F9 now a ratio, not a percent)
$ws.Range("F9").Formula = "=C9/C8"

# Add new labels describing what the ratios in column F are relative to
$ws.Range("G10").Value = "of FM"
$ws.Range("G9").Value = "of MB"

# Update the active selection to reflect where the user ended up after editing
$ws.Range("G11").Select()
